# Auto-generated Excel COM-interop script to apply the Ixion_Profits market-data refresh
# across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets (scheduled-runner "Update Sheets" commit).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 196
$ws.Range("I28").Value = 126.666664
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 126.666664
$ws.Range("L28").Value = 300
$ws.Range("M28").Value = 358.333336
$ws.Range("N28").Value = -1270
$ws.Range("H82").Value = 3500
$ws.Range("I82").Value = 3500
$ws.Range("K82").Value = 10500
$ws.Range("M82").Value = -10094
$ws.Range("H85").Value = 3500
$ws.Range("I85").Value = 3500
$ws.Range("K85").Value = 10500
$ws.Range("M85").Value = -9096
$ws.Range("H88").Value = 2849.5
$ws.Range("I88").Value = 2667.3333
$ws.Range("J88").Value = 2910.2222
$ws.Range("K88").Value = 2667.3333
$ws.Range("L88").Value = 2910.2222
$ws.Range("M88").Value = -2261.3333
$ws.Range("N88").Value = -3722.2222
$ws.Range("H91").Value = 2849.5
$ws.Range("I91").Value = 2667.3333
$ws.Range("J91").Value = 2910.2222
$ws.Range("K91").Value = 2667.3333
$ws.Range("L91").Value = 2910.2222
$ws.Range("M91").Value = -1263.3333
$ws.Range("N91").Value = -5718.2222
$ws.Range("H116").Value = 6832.8623
$ws.Range("I116").Value = 8756.412
$ws.Range("K116").Value = 8756.412
$ws.Range("M116").Value = -5314.412
$ws.Range("H138").Value = 2055.5474
$ws.Range("I138").Value = 1103.5714
$ws.Range("J138").Value = 2610.8667
$ws.Range("K138").Value = 3310.7142
$ws.Range("L138").Value = 7832.6001
$ws.Range("M138").Value = 1829.2858
$ws.Range("N138").Value = -18112.6001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1744.9474
$ws.Range("I2").Value = 1674.5
$ws.Range("K2").Value = 1674.5
$ws.Range("M2").Value = -1561.5
$ws.Range("H22").Value = 2110.5
$ws.Range("I22").Value = 2110.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2110.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1811.5
$ws.Range("N22").ClearContents()
$ws.Range("H88").Value = 2049.8333
$ws.Range("I88").Value = 1575
$ws.Range("K88").Value = 1575
$ws.Range("M88").Value = -1169
$ws.Range("H91").Value = 2049.8333
$ws.Range("I91").Value = 1575
$ws.Range("K91").Value = 1575
$ws.Range("M91").Value = -171
$ws.Range("H95").Value = 43000
$ws.Range("J95").Value = 43000
$ws.Range("L95").Value = 43000
$ws.Range("N95").Value = -48492
$ws.Range("H116").Value = 1744.9474
$ws.Range("I116").Value = 1674.5
$ws.Range("K116").Value = 1674.5
$ws.Range("M116").Value = 619.5
$ws.Range("H122").Value = 1604693.2
$ws.Range("I122").Value = 1604693.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4814079.6
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4811629.6
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1410702.8
$ws.Range("I132").Value = 1597.0944
$ws.Range("J132").Value = 5559736
$ws.Range("K132").Value = 4791.2832
$ws.Range("L132").Value = 16679208
$ws.Range("M132").Value = -2261.2832
$ws.Range("N132").Value = -16684268

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1744.9474
$ws.Range("I3").Value = 1674.5
$ws.Range("K3").Value = 1674.5
$ws.Range("M3").Value = -1560.5
$ws.Range("H124").Value = 20023.592
$ws.Range("I124").Value = 9905.263000000001
$ws.Range("K124").Value = 9905.263000000001
$ws.Range("M124").Value = -4995.263000000001
$ws.Range("H126").Value = 57846.25
$ws.Range("J126").Value = 57846.25
$ws.Range("L126").Value = 57846.25
$ws.Range("N126").Value = -67726.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3764
$ws.Range("I122").Value = 3764
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11292
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8842
$ws.Range("N122").ClearContents()
$ws.Range("H138").Value = 62000
$ws.Range("J138").Value = 62000
$ws.Range("L138").Value = 62000
$ws.Range("N138").Value = -72280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 270.2857
$ws.Range("I2").Value = 276.66666
$ws.Range("K2").Value = 276.66666
$ws.Range("M2").Value = -163.66666
$ws.Range("H35").Value = 250925
$ws.Range("I35").Value = 250925
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 250925
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -250627
$ws.Range("N35").ClearContents()
$ws.Range("H57").Value = 8142.4443
$ws.Range("J57").Value = 8142.4443
$ws.Range("L57").Value = 8142.4443
$ws.Range("N57").Value = -9782.444299999999
$ws.Range("H80").Value = 6917.5
$ws.Range("I80").Value = 8979
$ws.Range("K80").Value = 8979
$ws.Range("M80").Value = -7981
$ws.Range("H83").Value = 6917.5
$ws.Range("I83").Value = 8979
$ws.Range("K83").Value = 44895
$ws.Range("M83").Value = -39903
$ws.Range("H110").Value = 27142.857
$ws.Range("J110").Value = 27142.857
$ws.Range("L110").Value = 27142.857
$ws.Range("N110").Value = -35322.857
$ws.Range("H113").Value = 125002500
$ws.Range("I113").Value = 333334660
$ws.Range("J113").Value = 3200
$ws.Range("K113").Value = 333334660
$ws.Range("L113").Value = 3200
$ws.Range("M113").Value = -333332490
$ws.Range("N113").Value = -7540
$ws.Range("H122").Value = 88736710
$ws.Range("I122").Value = 88736710
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 266210130
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -266207680
$ws.Range("N122").ClearContents()
$ws.Range("H139").Value = 55214.133
$ws.Range("J139").Value = 55214.133
$ws.Range("L139").Value = 55214.133
$ws.Range("N139").Value = -65494.133

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2305.5334
$ws.Range("I7").Value = 1934.3636
$ws.Range("K7").Value = 1934.3636
$ws.Range("M7").Value = -1822.3636
$ws.Range("H16").Value = 1964.5714
$ws.Range("I16").Value = 1964.5714
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1964.5714
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1794.5714
$ws.Range("N16").ClearContents()
$ws.Range("H33").Value = 10000
$ws.Range("J33").Value = 10000
$ws.Range("L33").Value = 10000
$ws.Range("N33").Value = -10580
$ws.Range("H35").Value = 4362.3335
$ws.Range("I35").Value = 2060.3333
$ws.Range("J35").Value = 8966.333000000001
$ws.Range("K35").Value = 2060.3333
$ws.Range("L35").Value = 8966.333000000001
$ws.Range("M35").Value = -1724.3333
$ws.Range("N35").Value = -9638.333000000001
$ws.Range("H61").Value = 1451.8667
$ws.Range("I61").Value = 1451.8667
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1451.8667
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1249.8667
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 1451.8667
$ws.Range("I113").Value = 1451.8667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1451.8667
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 718.1333
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 3251892
$ws.Range("I122").Value = 3251892
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9755676
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9753226
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 2305.5334
$ws.Range("I126").Value = 1934.3636
$ws.Range("K126").Value = 5803.0908
$ws.Range("M126").Value = -3333.0908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2161.8708
$ws.Range("I113").Value = 2044.0416
$ws.Range("J113").Value = 2565.8572
$ws.Range("K113").Value = 6132.1248
$ws.Range("L113").Value = 7697.571599999999
$ws.Range("M113").Value = -3962.1248
$ws.Range("N113").Value = -12037.5716
$ws.Range("H116").Value = 44999.5
$ws.Range("J116").Value = 44999.5
$ws.Range("L116").Value = 44999.5
$ws.Range("N116").Value = -54177.5
$ws.Range("H122").Value = 1873.4445
$ws.Range("I122").Value = 1873.4445
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5620.333500000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3170.333500000001
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 775.2083
$ws.Range("I126").Value = 663.86365
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 1991.59095
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = 478.40905
$ws.Range("N126").Value = -10940
